$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.122.86'
$ws.Range("E2").Value = '  +6.29%  '
$ws.Range("D3").Value = '3.115.54'
$ws.Range("E3").Value = '  +3.93%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("E5").Value = '  +4.53%  '
$ws.Range("D6").Value = '144.35'
$ws.Range("E6").Value = '  +4.07%  '
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("D8").Value = '3.106.03'
$ws.Range("E8").Value = '  +4.15%  '
$ws.Range("E9").Value = '  +2.36%  '
$ws.Range("E10").Value = '  +10.90%  '
$ws.Range("D11").Value = '''5.70'
$ws.Range("E11").Value = '  +10.24%  '
$ws.Range("D12").Value = '''0.470'
$ws.Range("E12").Value = '  +2.27%  '
$ws.Range("D13").Value = '0.0000246'
$ws.Range("E13").Value = '  +6.77%  '
$ws.Range("D14").Value = '35.61'
$ws.Range("E14").Value = '  +5.51%  '
$ws.Range("D15").Value = '0.123'
$ws.Range("E15").Value = '  +0.77%  '
$ws.Range("D16").Value = '3.629.73'
$ws.Range("E16").Value = '  +3.92%  '
$ws.Range("D17").Value = '7.27'
$ws.Range("E17").Value = '  -0.49%  '
$ws.Range("D18").Value = '63.056.45'
$ws.Range("E18").Value = '  +6.21%  '
$ws.Range("D19").Value = '3.110.47'
$ws.Range("E19").Value = '  +3.81%  '
$ws.Range("D20").Value = '455.64'
$ws.Range("E20").Value = '  +5.73%  '
$ws.Range("D21").Value = '14.15'
$ws.Range("E21").Value = '  +3.95%  '
$ws.Range("D22").Value = '0.735'
$ws.Range("E22").Value = '  +1.90%  '
$ws.Range("E23").Value = '  +6.60%  '
$ws.Range("D24").Value = '13.66'
$ws.Range("E24").Value = '  +0.89%  '
$ws.Range("D25").Value = '82.32'
$ws.Range("E25").Value = '  +2.42%  '
$ws.Range("E26").Value = '  +0.13%  '
$ws.Range("E27").Value = '  +1.51%  '
$ws.Range("E28").Value = '  +6.51%  '
$ws.Range("E29").Value = '  +5.48%  '
$ws.Range("E30").Value = '  +0.01%  '
$ws.Range("E31").Value = '  +12.53%  '
$ws.Range("E32").Value = '  +10.98%  '
$ws.Range("E33").Value = '  +5.25%  '
$ws.Range("B34").Value = 'Stacks'
$ws.Range("C34").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D34").Value = '2.37'
$ws.Range("E34").Value = '  +12.35%  '
$ws.Range("B35").Value = 'PEPE'
$ws.Range("C35").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D35").Value = '0.0₃0823'
$ws.Range("E35").Value = '  +7.61%  '
$ws.Range("E36").Value = '  +3.37%  '
$ws.Range("D37").Value = '6.08'
$ws.Range("E37").Value = '  +1.94%  '
$ws.Range("D38").Value = '3.15'
$ws.Range("E38").Value = '  +14.46%  '
$ws.Range("D39").Value = '51.04'
$ws.Range("E39").Value = '  +4.45%  '
$ws.Range("E40").Value = '  +1.28%  '
$ws.Range("D41").Value = '429.74'
$ws.Range("E41").Value = '  +5.62%  '
$ws.Range("D42").Value = '2.972.50'
$ws.Range("E42").Value = '  +7.08%  '
$ws.Range("D43").Value = '0.0375'
$ws.Range("E43").Value = '  +5.87%  '
$ws.Range("E44").Value = '  +4.52%  '
$ws.Range("D45").Value = '0.276'
$ws.Range("E45").Value = '  +9.41%  '
$ws.Range("D46").Value = '2.18'
$ws.Range("E46").Value = '  +8.50%  '
$ws.Range("D47").Value = '125.11'
$ws.Range("E47").Value = '  +1.40%  '
$ws.Range("D49").Value = '34.72'
$ws.Range("E49").Value = '  -0.21%  '
$ws.Range("E50").Value = '  +1.33%  '
$ws.Range("D51").Value = '24.95'
$ws.Range("E51").Value = '  +5.91%  '
